$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting existing rows 22-38 down to 23-39
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new data record
$ws.Cells.Item(22, 1).Value = 7
$ws.Cells.Item(22, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(22, 3).Value = "Ñuble"
$ws.Cells.Item(22, 4).Value = 44762
$ws.Cells.Item(22, 5).Value = 16
$ws.Cells.Item(22, 6).Value = 100112013
$ws.Cells.Item(22, 7).Value = "Alcachofa"
$ws.Cells.Item(22, 8).Value = "Argentina(o)"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 60
$ws.Cells.Item(22, 11).Value = 15000
$ws.Cells.Item(22, 12).Value = 16000
$ws.Cells.Item(22, 13).Value = 15500
$ws.Cells.Item(22, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(22, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(22, 16).Value = 310
$ws.Cells.Item(22, 17).Value = 50
$ws.Cells.Item(22, 18).Value = "Hortaliza"
